$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.543.78'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.510.18'
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.79'
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.91'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("E7").Value = '  +2.21%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("E9").Value = '  -1.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.46'
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.73'
$ws.Range("E12").Value = '  +3.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.113'
$ws.Range("E13").Value = '  -3.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.895.55'
$ws.Range("E14").Value = '  -1.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.67'
$ws.Range("E15").Value = '  +8.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.501.22'
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.860'
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.533.89'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.87'
$ws.Range("E19").Value = '  -4.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0972'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.48'
$ws.Range("E21").Value = '  -1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.46'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '253.08'
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.95'
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("E25").Value = '  -2.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.99'
$ws.Range("E26").Value = '  -3.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  +10.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.18'
$ws.Range("E29").Value = '  +1.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.60'
$ws.Range("E30").Value = '  -3.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.92'
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.42'
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.20'
$ws.Range("E33").Value = '  +4.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.29'
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.07'
$ws.Range("E35").Value = '  -4.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0786'
$ws.Range("E36").Value = '  -2.05%  '
$ws.Range("E37").Value = '  -4.76%  '
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.77'
$ws.Range("E39").Value = '  -4.04%  '
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.87'
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.38'
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.02'
$ws.Range("E43").Value = '  -1.66%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0302'
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.027.62'
$ws.Range("E46").Value = '  -1.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.29'
$ws.Range("E47").Value = '  -4.58%  '
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.754.52'
$ws.Range("E49").Value = '  -1.82%  '
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.66'
$ws.Range("E51").Value = '  -4.55%  '
